$wb = $excel.ActiveWorkbook

# --- Add the new "Player Info" sheet (inserted before the current first/active sheet) ---
$new = $wb.Worksheets.Add()
$new.Name = "Player Info"

$new.Range("A1").Value = "ID"
$new.Range("B1").Value = "NAME"
$new.Range("C1").Value = "BATTING_HAND"
$new.Range("D1").Value = "BOWL_STYLE"

$new.Range("A2").NumberFormat = "@"
$new.Range("A2").Value = "4611"
$new.Range("A2").ClearFormats()
$new.Range("B2").Value = "Karim Janat"
$new.Range("C2").Value = "Right Handed"
$new.Range("D2").Value = "Right Arm Medium"

# copy the header formatting (bold, border, alignment) from the "ODI Batting" sheet
$wb.Worksheets.Item("ODI Batting").Range("A1:D1").Copy()
$new.Range("A1:D1").PasteSpecial(-4122)

# --- Rename MATCH_CARD_LINK -> MATCH_CODE and replace the URL value with the bare match code ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "3996"
$batting.Range("D2").ClearFormats()

$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "3996"
$bowling.Range("B2").ClearFormats()
